$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errMsg = " Erreur : 502 Server Error: Bad Gateway for url: https://chatbot-o4gm.onrender.com/generate_email"

$ws.Range("A3").Value = "CIC"
$ws.Range("B3").Value = $errMsg

$ws.Range("A4").Value = "Claire Huteau"
$ws.Range("B4").Value = $errMsg
